$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Implemented Features")

# Rename "Implemented Features" -> "Features  To dos"
$ws2.Name = "Features  To dos"

# Add a new changelog row (row 10) documenting the mean/var addition
$ws2.Cells.Item(10, 1).Value = "f.calc_variability_seg_M_joint"
$ws2.Cells.Item(10, 2).Value = 44643
$ws2.Cells.Item(10, 3).Value = "Thomas"
$ws2.Cells.Item(10, 4).Value = "mean / var if data set sufficent"
$ws2.Cells.Item(10, 5).Value = 44643

# Header D1: "Comments" -> "Changes"
$ws2.Range("D1").Value = "Changes"

# Widen column D to fit the new, longer "Changes" text
$ws2.Columns.Item(4).ColumnWidth = 50.9

$ws2.Cells.Item(10, 6).Value = "main"

# Copy the date number formats from row 9 so the new row matches existing styling
$ws2.Cells.Item(9, 2).Copy()
$ws2.Cells.Item(10, 2).PasteSpecial(-4122)
$ws2.Cells.Item(9, 5).Copy()
$ws2.Cells.Item(10, 5).PasteSpecial(-4122)

# Update selection on the "Features  To dos" sheet
$ws2.Range("E11").Select()

# Make "Features  To dos" the active/visible sheet (moves tabSelected + activeTab)
$ws2.Activate()
